$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = 6.096
$ws.Range("B12").Value = 4.951000000000001
$ws.Range("B18").Value = 5.194
$ws.Range("B37").Value = 8.73
$ws.Range("B55").Value = 4.763
$ws.Range("B68").Value = 5.139
$ws.Range("B77").Value = 5.614
$ws.Range("B78").Value = 7.711
